$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.309.12'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.840.95'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '363.62'
$ws.Range("E5").Value = '  +3.51%  '
$ws.Range("D6").Value = '112.25'
$ws.Range("E6").Value = '  -3.15%  '
$ws.Range("E7").Value = '  +3.72%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.604'
$ws.Range("E9").Value = '  +2.33%  '
$ws.Range("D10").Value = '41.01'
$ws.Range("E10").Value = '  -3.48%  '
$ws.Range("D11").Value = '0.0867'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  +1.05%  '
$ws.Range("D13").Value = '20.03'
$ws.Range("E14").Value = '  -0.71%  '
$ws.Range("D15").Value = '3.279.77'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '2.876.57'
$ws.Range("E16").Value = '  +2.90%  '
$ws.Range("E17").Value = '  +4.63%  '
$ws.Range("D18").Value = '52.158.24'
$ws.Range("E18").Value = '  +0.00%  '
$ws.Range("D19").Value = '7.50'
$ws.Range("E19").Value = '  +2.87%  '
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").Value = '0.0000100'
$ws.Range("E22").Value = '  +2.07%  '
$ws.Range("D23").Value = '273.24'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").Value = '70.52'
$ws.Range("E24").Value = '  +0.37%  '
$ws.Range("D25").Value = '2.84'
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("D26").Value = '26.99'
$ws.Range("E26").Value = '  +0.58%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = '10.33'
$ws.Range("E28").Value = '  +0.74%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("E31").Value = '  +3.25%  '
$ws.Range("B32").Value = 'InjectiveProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D32").Value = '35.23'
$ws.Range("E32").Value = '  +2.27%  '
$ws.Range("B33").Value = 'OKB'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D33").Value = '52.39'
$ws.Range("E33").Value = '  +4.26%  '
$ws.Range("D34").Value = '5.89'
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("D35").Value = '5.60'
$ws.Range("E35").Value = '  +13.12%  '
$ws.Range("D36").Value = '0.0854'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  +1.69%  '
$ws.Range("D39").Value = '2.05'
$ws.Range("E39").Value = '  -2.85%  '
$ws.Range("D40").Value = '18.43'
$ws.Range("E40").Value = '  -1.74%  '
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").Value = '2.55'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '125.15'
$ws.Range("E43").Value = '  -0.84%  '
$ws.Range("D44").Value = '22.79'
$ws.Range("E44").Value = '  -2.60%  '
$ws.Range("E45").Value = '  -1.23%  '
$ws.Range("D46").Value = '2.086.00'
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("D47").Value = '3.34'
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = '2.33'
$ws.Range("E48").Value = '  +1.48%  '
$ws.Range("D49").Value = '5.89'
$ws.Range("E49").Value = '  +5.62%  '
$ws.Range("D50").Value = '0.974'
$ws.Range("E50").Value = '  +2.46%  '
$ws.Range("D51").Value = '9.21'
$ws.Range("E51").Value = '  +2.50%  '
